$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.390.30'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '2.286.88'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '300.99'
$ws.Range("E5").Value = '  -1.99%  '

$ws.Range("D6").Value = "'95.30"
$ws.Range("E6").Value = '  -1.33%  '

$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  -2.52%  '

$ws.Range("D10").Value = '34.36'
$ws.Range("E10").Value = '  -3.08%  '

$ws.Range("D11").Value = '18.99'
$ws.Range("E11").Value = '  +2.51%  '

$ws.Range("D12").Value = '0.0779'
$ws.Range("E12").Value = '  -1.65%  '

$ws.Range("D13").Value = '0.118'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("D14").Value = '6.72'
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("D15").Value = '2.645.92'
$ws.Range("E15").Value = '  -0.63%  '

$ws.Range("D16").Value = '2.290.30'
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").Value = '0.777'
$ws.Range("E17").Value = '  -0.85%  '

$ws.Range("D18").Value = '42.348.67'
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("D19").Value = '12.14'
$ws.Range("E19").Value = '  -7.64%  '

$ws.Range("D20").Value = '0.0₃0886'
$ws.Range("E20").Value = '  -1.42%  '

$ws.Range("D21").Value = '5.95'
$ws.Range("E21").Value = '  -1.63%  '

$ws.Range("D22").Value = '67.64'
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").Value = '2.26'
$ws.Range("E23").Value = '  +5.24%  '

$ws.Range("D24").Value = '234.94'

$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").Value = "'2.40"

$ws.Range("D27").Value = '24.24'
$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("E28").Value = '  -1.12%  '

$ws.Range("D29").Value = '165.17'
$ws.Range("E29").Value = '  -0.66%  '

$ws.Range("D30").Value = '9.02'
$ws.Range("E30").Value = '  -0.70%  '

$ws.Range("D31").Value = '31.79'
$ws.Range("E31").Value = '  -4.36%  '

$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").Value = '4.97'
$ws.Range("E33").Value = '  -0.59%  '

$ws.Range("D34").Value = '17.48'
$ws.Range("E34").Value = '  -1.87%  '

$ws.Range("D35").Value = '0.0694'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.32'
$ws.Range("E36").Value = '  -2.85%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.36'
$ws.Range("E37").Value = '  -8.34%  '

$ws.Range("D38").Value = '0.0996'
$ws.Range("E38").Value = '  -1.55%  '

$ws.Range("D39").Value = '1.73'
$ws.Range("E39").Value = '  -1.70%  '

$ws.Range("E40").Value = '  -1.56%  '

$ws.Range("D41").Value = '2.67'
$ws.Range("E41").Value = '  -1.60%  '

$ws.Range("D42").Value = '19.88'
$ws.Range("E42").Value = '  +8.64%  '

$ws.Range("D43").Value = '1.958.55'

$ws.Range("D44").Value = '10.33'
$ws.Range("E44").Value = '  +2.78%  '

$ws.Range("D45").Value = '0.0277'
$ws.Range("E45").Value = '  -0.93%  '

$ws.Range("E46").Value = '  +0.35%  '

$ws.Range("D47").Value = '2.73'
$ws.Range("E47").Value = '  -2.87%  '

$ws.Range("D48").Value = '2.516.13'
$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("D49").Value = '52.89'
$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("D50").Value = '2.79'
$ws.Range("E50").Value = '  -3.07%  '

$ws.Range("D51").Value = '70.71'
$ws.Range("E51").Value = '  -1.28%  '
